$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I (shifts existing I..N to J..O)
$ws.Columns("I").Insert()

# Set the new column header and values
$ws.Range("I1").Value = "Onboarding Completed"
$ws.Range("I2").Value = "Yes"
$ws.Range("I3").Value = "No"
$ws.Range("I4").Value = "Yes"
$ws.Range("I5").Value = "Yes"
$ws.Range("I6").Value = "Yes"
$ws.Range("I7").Value = "Yes"

# Clear the redundant style that had been applied to K1:N1 (now shifted to L1:O1) (now default/Normal)
$ws.Range("L1:O1").Style = "Normal"

# Update the selection to reflect the new onboarding column
$ws.Range("I1:I7").Select()
